$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Data de entrega:" 06/07 -> 14/12
# ------------------------------------------------------------------
$d.Content.Find.Execute("06/07", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "14/12", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "...enviada para o monitor da disciplina." ->
#    "...enviada para o professor da disciplina."
#    (the replaced word becomes its own run, and the _GoBack bookmark
#    ends up sitting right after it, exactly like Word leaves it after
#    a user types a replacement over a selection)
# ------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("monitor")
$len = "monitor".Length
$r = $d.Range($idx, $idx + $len)
$r.Text = "professor"

$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("professor da disciplina")
$posBeforeProfessor = $idx2
$posAfterProfessor = $idx2 + "professor".Length

# Force a run boundary before "professor" with a scratch bookmark,
# then drop the scratch bookmark once the split exists.
$bmBeforeRange = $d.Range($posBeforeProfessor, $posBeforeProfessor)
$d.Bookmarks.Add("ZZZScratchSplit", $bmBeforeRange) | Out-Null

# Re-seat _GoBack right after "professor" (this also forces the run
# boundary right after the word, like Word does for the last edit).
$bmAfterRange = $d.Range($posAfterProfessor, $posAfterProfessor)
$d.Bookmarks.Add("_GoBack", $bmAfterRange) | Out-Null

# Remove the scratch bookmark now that the split it created persists.
$d.Bookmarks("ZZZScratchSplit").Delete()
